$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell B20: fix apostrophe case "Developer'S" -> "Developer's"
$ws.Cells.Item(20, 2).Value = "Unit And Integration Tests. Numerical Analysis And Gradient Checking. Research Papers And Its Source Codes. Consultation With Experts. Uses Developer's University Lecturer Notes."

# Update cell A9: "Minimum Lines To Train A Model" -> "Minimum Number Of Lines To Train A Model"
$ws.Cells.Item(9, 1).Value = "Minimum Number Of Lines To Train A Model"

# Update the selected/active cell to B12 (as reflected in the sheet view state)
$ws.Range("B12").Select()
